# Generate Report for Archive
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn/de-de status columns) and on each language
#   sheet's Status column.
# - Narrow the now-shorter status columns to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Resize the status columns to fit the new, shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.576851254417766
$overview.Columns.Item(6).ColumnWidth = 12.576851254417766

$zhcn.Columns.Item(3).ColumnWidth = 12.576851254417766
$dede.Columns.Item(3).ColumnWidth = 12.576851254417766
